# Refresh the "cryptos" symbol list (GitHub Actions scrape) — prices, the
# 24h-ranked coin/link/volume rows shift by one slot, and the capture hour
# moves from 12 -> 13 for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking text ("248.92", "13",
# "0.1400" with a significant trailing zero, etc.). Mark them as Text first
# so the Value setter below stores them verbatim instead of coercing to
# Double (which would both change the type and drop trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '248.92'
$ws.Range("G2").Value = '13'

$ws.Range("D3").Value = '22.71'
$ws.Range("G3").Value = '13'

$ws.Range("D4").Value = '5.281'
$ws.Range("G4").Value = '13'

$ws.Range("D5").Value = '0.05685'
$ws.Range("G5").Value = '13'

$ws.Range("G6").Value = '13'

$ws.Range("D7").Value = '6.346'
$ws.Range("G7").Value = '13'

$ws.Range("D8").Value = '0.8062'
$ws.Range("G8").Value = '13'

$ws.Range("D9").Value = '0.9041'
$ws.Range("G9").Value = '13'

$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D10").Value = '0.0005812'
$ws.Range("E10").Value = '9OneONE'
$ws.Range("G10").Value = '13'

$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1400'
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("G11").Value = '13'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.07453'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("G12").Value = '13'

$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").Value = '0.03090'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("G13").Value = '13'

$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Value = '0.03011'
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("G14").Value = '13'

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = '0.09374'
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("G15").Value = '13'

$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").Value = '3.863'
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("G16").Value = '13'

$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").Value = '0.001574'
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("G17").Value = '13'

$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").Value = '0.04782'
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("G18").Value = '13'

$ws.Range("B19").Value = 'UpBots'
$ws.Range("C19").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D19").Value = '0.01829'
$ws.Range("E19").Value = '18UpBotsUBXTBestin24h'
$ws.Range("G19").Value = '13'

$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").Value = '0.006437'
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("G20").Value = '13'

$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D21").Value = '0.004986'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("G21").Value = '13'

$ws.Range("B22").Value = 'BitKan'
$ws.Range("C22").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D22").Value = '0.001000'
$ws.Range("E22").Value = '21BitKanKAN'
$ws.Range("G22").Value = '13'

$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").Value = '0.0001500'
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("G23").Value = '13'

$ws.Range("B24").Value = 'LEO'
$ws.Range("C24").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D24").Value = '3.694'
$ws.Range("E24").Value = '23LEOLEO'
$ws.Range("G24").Value = '13'

$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").Value = '2.201'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("G25").Value = '13'

$ws.Range("D26").Value = '0.3257'
$ws.Range("G26").Value = '13'

$ws.Range("G27").Value = '13'

$ws.Range("G28").Value = '13'

$ws.Range("G29").Value = '13'

$ws.Range("G30").Value = '13'

$ws.Range("G31").Value = '13'

$ws.Range("G32").Value = '13'

$ws.Range("G33").Value = '13'

$ws.Range("G34").Value = '13'

$ws.Range("G35").Value = '13'

$ws.Range("G36").Value = '13'

$ws.Range("G37").Value = '13'

$ws.Range("G38").Value = '13'

$ws.Range("G39").Value = '13'

$ws.Range("D40").Value = '0.03980'
$ws.Range("G40").Value = '13'

$ws.Range("D41").Value = '0.006721'
$ws.Range("G41").Value = '13'

$ws.Range("G42").Value = '13'

$ws.Range("D43").Value = '0.002765'
$ws.Range("G43").Value = '13'

$ws.Range("D44").Value = '0.008128'
$ws.Range("G44").Value = '13'

$ws.Range("D45").Value = '0.00005581'
$ws.Range("G45").Value = '13'

$ws.Range("G46").Value = '13'

$ws.Range("G47").Value = '13'

$ws.Range("D48").Value = '0.2042'
$ws.Range("G48").Value = '13'

$ws.Range("G49").Value = '13'

$ws.Range("G50").Value = '13'

$ws.Range("G51").Value = '13'
